# Add separate analysis rows for elderly / non-elderly lethality (rows 20-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing formatted cell A19 as the format source for the new
# column-A "index" cells so the new rows reuse the same cell style
# (bold font, thin box border, centered/top alignment) instead of
# creating a brand new style entry.
$formatSource = $ws.Range("A19")

$rows = @(
    @{ Row = 20; A = 18; B = "Lethality_eld_tx";              C = 17.32805314924295;  D = 14.96524615730777;  E = 19.73942139832099;  F = 16862; G = 22215; H = 23033; I = 36.5970821966552 },
    @{ Row = 21; A = 19; B = "Lethality_non_eld_tx";          C = 15.43063821181252;  D = 12.87963735978455;  E = 18.03928989882952;  F = 6834;  G = 8212;  H = 9097;  I = 33.11384255194616 },
    @{ Row = 22; A = 20; B = "Lethality_uti_eld_tx";          C = 0.1016167906223275; D = -4.147032013637064; E = 4.538585446015642;  F = 57030; G = 56546; H = 57160; I = 0.2279502016482553 },
    @{ Row = 23; A = 21; B = "Lethality_uti_non_eld_tx";      C = -0.5005529586716739;D = -4.144797049125626; E = 3.282238801412829;  F = 36040; G = 36190; H = 35675; I = -1.01276359600444 },
    @{ Row = 24; A = 22; B = "Lethality_non_uti_eld_tx";      C = 17.42319019148113;  D = 14.69182539274845;  E = 20.21960194223688;  F = 13754; G = 17989; H = 18804; I = 36.71659153700742 },
    @{ Row = 25; A = 23; B = "Lethality_non_uti_non_eld_tx";  C = 13.39894678498341;  D = 10.2772578220238;   E = 16.60900339666704;  F = 4926;  G = 5581;  H = 6333;  I = 28.56272838002436 }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Column A: copy the formatting from the existing styled index column,
    # then write the value.
    $formatSource.Copy()
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowIndex, 1).Value = $r.A

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    $ws.Cells.Item($rowIndex, 5).Value = $r.E
    $ws.Cells.Item($rowIndex, 6).Value = $r.F
    $ws.Cells.Item($rowIndex, 7).Value = $r.G
    $ws.Cells.Item($rowIndex, 8).Value = $r.H
    $ws.Cells.Item($rowIndex, 9).Value = $r.I
}

$excel.CutCopyMode = 0
